$d = $word.ActiveDocument

# The paragraph we need to touch is the one that carries the "_GoBack"
# bookmark (it sits exactly at the boundary between the two runs
# "There will be a conflict in" and " this file").
$bm = $d.Bookmarks.Item("_GoBack")
$boundary = $bm.Start

# Find the paragraph that contains that boundary.
$paraIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Start -le $boundary -and $boundary -le $p.Range.End) {
        $paraIndex = $i
    }
}
$para = $d.Paragraphs($paraIndex)

# Range covering the first run ("There will be a conflict in") and the
# range covering the second run (" this file"), located purely from the
# bookmark position + the paragraph boundaries (no hard-coded offsets).
$run1 = $d.Range($para.Range.Start, $boundary)
$run2Len = $para.Range.End - $boundary - 1   # -1 excludes the paragraph mark
$run2 = $d.Range($boundary, $boundary + $run2Len)

$run1Text = $run1.Text
$run2Text = $run2.Text

# 1) "There will be a conflict in" + " this file" -> merge into a single
#    run so paragraph A reads "There will be a conflict in this file".
$run1.Text = $run1Text + $run2Text
$newBoundary = $run1.End

# 2) The old second-run span (now shifted right after the merge) still
#    holds the stale " this file" text - turn it into "ng line".
$tail = $d.Range($newBoundary, $newBoundary + $run2Len)
$tail.Text = "ng line"

# 3) Insert the new leading text "Added conflicti" right before the
#    bookmark boundary; it merges into the preceding run and therefore
#    inherits its run formatting (language etc.).
$insPoint = $d.Range($newBoundary, $newBoundary)
$insPoint.InsertBefore("Added conflicti")

# 4) Finally split the paragraph at that same boundary so "Added
#    conflicti<bookmark>ng line" becomes its own new paragraph, right
#    after "There will be a conflict in this file".
$splitPoint = $d.Range($newBoundary, $newBoundary)
$splitPoint.InsertParagraphAfter()

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
foreach ($p in $d.Paragraphs) {
    Write-Output $p.Range.Text
}
